$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl20"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07886166666666666
$ws.Range("H2").Value = 0.236585
$ws.Range("I2").Value = 0.5304501264551309
$ws.Range("J2").Value = 0.5304501264551309
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1143813333333333
$ws.Range("N2").Value = 0.343144
$ws.Range("O2").Value = 0.128300337591142
$ws.Range("P2").Value = 0.1283003375911419
$ws.Range("Q2").Value = 0.009020302582222223
$ws.Range("R2").Value = 0.08118272324
$ws.Range("S2").Value = 0.06805693029945724
$ws.Range("T2").Value = 0.06805693029945722

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl20"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07886166666666666
$ws.Range("H3").Value = 0.236585
$ws.Range("I3").Value = 0.5304501264551309
$ws.Range("J3").Value = 0.5304501264551309
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7200953333333334
$ws.Range("N3").Value = 2.160286
$ws.Range("O3").Value = 0.8077233554817153
$ws.Range("P3").Value = 0.8077233554817151
$ws.Range("Q3").Value = 0.05678791814555556
$ws.Range("R3").Value = 0.51109126331
$ws.Range("S3").Value = 0.4284569560560385
$ws.Range("T3").Value = 0.4284569560560385

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl20"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07886166666666666
$ws.Range("H4").Value = 0.236585
$ws.Range("I4").Value = 0.5304501264551309
$ws.Range("J4").Value = 0.5304501264551309
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03357866666666667
$ws.Range("N4").Value = 0.100736
$ws.Range("O4").Value = 0.03766483694187069
$ws.Range("P4").Value = 0.03766483694187069
$ws.Range("Q4").Value = 0.002648069617777778
$ws.Range("R4").Value = 0.02383262656
$ws.Range("S4").Value = 0.01997931751872719
$ws.Range("T4").Value = 0.01997931751872719

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ccl20"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07886166666666666
$ws.Range("H5").Value = 0.236585
$ws.Range("I5").Value = 0.5304501264551309
$ws.Range("J5").Value = 0.5304501264551309
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023457
$ws.Range("N5").Value = 0.070371
$ws.Range("O5").Value = 0.02631146998527222
$ws.Range("P5").Value = 0.02631146998527222
$ws.Range("Q5").Value = 0.001849858115
$ws.Range("R5").Value = 0.016648723035
$ws.Range("S5").Value = 0.01395692258090803
$ws.Range("T5").Value = 0.01395692258090803

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ccl20"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.06980766666666667
$ws.Range("H6").Value = 0.209423
$ws.Range("I6").Value = 0.4695498735448692
$ws.Range("J6").Value = 0.4695498735448692
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1143813333333333
$ws.Range("N6").Value = 0.343144
$ws.Range("O6").Value = 0.128300337591142
$ws.Range("P6").Value = 0.1283003375911419
$ws.Range("Q6").Value = 0.007984693990222223
$ws.Range("R6").Value = 0.07186224591199999
$ws.Range("S6").Value = 0.06024340729168474
$ws.Range("T6").Value = 0.06024340729168472

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ccl20"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.06980766666666667
$ws.Range("H7").Value = 0.209423
$ws.Range("I7").Value = 0.4695498735448692
$ws.Range("J7").Value = 0.4695498735448692
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7200953333333334
$ws.Range("N7").Value = 2.160286
$ws.Range("O7").Value = 0.8077233554817153
$ws.Range("P7").Value = 0.8077233554817151
$ws.Range("Q7").Value = 0.05026817499755556
$ws.Range("R7").Value = 0.4524135749780001
$ws.Range("S7").Value = 0.3792663994256769
$ws.Range("T7").Value = 0.3792663994256768

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ccl20"
$ws.Range("C8").Value = "Ackr4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.06980766666666667
$ws.Range("H8").Value = 0.209423
$ws.Range("I8").Value = 0.4695498735448692
$ws.Range("J8").Value = 0.4695498735448692
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03357866666666667
$ws.Range("N8").Value = 0.100736
$ws.Range("O8").Value = 0.03766483694187069
$ws.Range("P8").Value = 0.03766483694187069
$ws.Range("Q8").Value = 0.002344048369777778
$ws.Range("R8").Value = 0.021096435328
$ws.Range("S8").Value = 0.0176855194231435
$ws.Range("T8").Value = 0.0176855194231435

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ccl20"
$ws.Range("C9").Value = "Ackr4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.06980766666666667
$ws.Range("H9").Value = 0.209423
$ws.Range("I9").Value = 0.4695498735448692
$ws.Range("J9").Value = 0.4695498735448692
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.023457
$ws.Range("N9").Value = 0.070371
$ws.Range("O9").Value = 0.02631146998527222
$ws.Range("P9").Value = 0.02631146998527222
$ws.Range("Q9").Value = 0.001637478437
$ws.Range("R9").Value = 0.014737305933
$ws.Range("S9").Value = 0.01235454740436419
$ws.Range("T9").Value = 0.01235454740436419
